# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Re-sort the "Periodo Mora" rows (17-20) into ascending order (1801..1804),
# keeping each period's own "Valor Mora" value attached to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E17").Value = "1801"
$ws.Range("F17").Value = 16666

$ws.Range("E18").Value = "1802"
$ws.Range("F18").Value = 31249

$ws.Range("E19").Value = "1803"
$ws.Range("F19").Value = 31249

$ws.Range("E20").Value = "1804"
$ws.Range("F20").Value = 31249
